$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows before row 308, shifting existing rows 308-413 down to 310-415.
$ws.Rows("308:309").Insert()

# Populate new row 308 (Primera, 2023-03-03)
$ws.Range("A308").Value = 7
$ws.Range("B308").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C308").Value = "Ñuble"
$ws.Range("D308").Value = 44988
$ws.Range("E308").Value = 16
$ws.Range("F308").Value = 100112023
$ws.Range("G308").Value = "Brócoli"
$ws.Range("H308").Value = "Sin especificar"
$ws.Range("I308").Value = "Primera"
$ws.Range("J308").Value = 300
$ws.Range("K308").Value = 900
$ws.Range("L308").Value = 900
$ws.Range("M308").Value = 900
$ws.Range("N308").Value = "$/unidad"
$ws.Range("O308").Value = "Región del Maule"
$ws.Range("P308").Value = 900
$ws.Range("Q308").Value = 1
$ws.Range("R308").Value = "Hortaliza"

# Populate new row 309 (Segunda, 2023-03-03)
$ws.Range("A309").Value = 7
$ws.Range("B309").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C309").Value = "Ñuble"
$ws.Range("D309").Value = 44988
$ws.Range("E309").Value = 16
$ws.Range("F309").Value = 100112023
$ws.Range("G309").Value = "Brócoli"
$ws.Range("H309").Value = "Sin especificar"
$ws.Range("I309").Value = "Segunda"
$ws.Range("J309").Value = 300
$ws.Range("K309").Value = 700
$ws.Range("L309").Value = 700
$ws.Range("M309").Value = 700
$ws.Range("N309").Value = "$/unidad"
$ws.Range("O309").Value = "Región del Maule"
$ws.Range("P309").Value = 700
$ws.Range("Q309").Value = 1
$ws.Range("R309").Value = "Hortaliza"
